$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @([double]"2.473e-19", 1.602, 0.526, 5),
    @([double]"1.81e-18", 2.529, 0.5928, 5),
    @(0.3462, 1.323, 0.2749, 20),
    @(0.3308, 0.74, 0.1534, 25),
    @(0.1351, 0.8452, 0.1706, 25),
    @(0.06888, 1.4, 0.4627, 20),
    @(0.05829, 1.679, 1.399, 15),
    @(0.2477, 2, 1.228, 20),
    @(0.4646, 1.501, 0.5923, 25),
    @(0.4228, 1.605, 0.1071, 25),
    @(0.3204, 1.748, 0.7953, 25),
    @(0.002187, 1.837, 0.1693, 10),
    @(0.1009, 2.557, 0.3006, 15),
    @(0.09429999999999999, 2.735, 0.242, 15),
    @(0.00134, 2.518, 0.7192, 10),
    @(0.0002662, 2.7, 0.2886, 10),
    @(0.007861, 3.241, 0.2737, 15),
    @(0.1738, 2.823, 0.2202, 25),
    @(0.01123, 3.415, 0.2173, 15),
    @(0.0126, 3.535, 0.107, 20),
    @(0.008581999999999999, 3.501, 0.1225, 20),
    @(0.03624, 3.327, 0.4257, 25),
    @(0.01635, 3.508, 0.6127, 25),
    @(0.009377, 4.187, 0.6296, 20),
    @(0.01314, 4.131, 0.1966, 20),
    @(0.01183, 4.956, 0.3537, 15),
    @(0.0006742, 4.715, 0.526, 15),
    @(0.02656, 4.982, 0.1488, 20),
    @(0.0008779, 5.018, 0.2842, 15),
    @(0.06807000000000001, 4.629, 0.253, 25),
    @(0.0004758, 4.611, 0.6179, 10),
    @(0.01389, 5.115, 0.2695, 15),
    @(0.005308, 5.212, 0.3045, 15),
    @(0.001326, 5.081, 0.164, 15),
    @(0.0005083, 5.216, 0.2042, 15),
    @(0.07881000000000001, 4.952, 0.2847, 25),
    @(0.06856, 4.969, 0.4637, 25),
    @(0.03904, 5.004, 0.7395, 25),
    @(0.01848, 5.541, 0.2089, 15),
    @(0.006495, 5.604, 0.08316999999999999, 15),
    @(0.003861, 5.56, 0.3786, 15),
    @(0.1434, 5.145, 0.2819, 25),
    @(0.229, 5.203, 1.989, 25),
    @(0.01825, 5.361, 2.31, 10),
    @(0.1221, 5.839, 0.4614, 15),
    @(0.2315, 5.903, 0.6021, 20),
    @(0.406, 5.547, 0.144, 25),
    @([double]"1.754e-18", 7.939, 0.5756, 5),
    @(0.03743, 6.44, 1.615, 15),
    @([double]"1.261e-17", 7.616, 0.9634, 5),
    @(0.02474, 6.173, 1.329, 15),
    @(0.08151, 5.96, 0.3352, 20),
    @(0.6843, 5.491, 0.4235, 25),
    @(0.3566, 5.392, 1.096, 25),
    @(0.07369000000000001, 5.366, 0.1319, 15),
    @(0.02181, 5.205, 0.6167, 15),
    @(0.01284, 4.392, 0.361, 10),
    @(0.243, 4.955, 0.3571, 25),
    @(0.02578, 4.763, 0.07553, 15),
    @(0.194, 4.809, 0.1818, 25),
    @(0.147, 4.719, 0.1987, 25),
    @(0.1728, 4.714, 0.9268999999999999, 25),
    @(0.006758, 4.34, 0.5497, 10),
    @(0.009176, 4.839, 0.2289, 15),
    @(0.002954, 4.865, 1.317, 15),
    @([double]"9.974e-05", 5.142, 0.3696, 10),
    @(0.1508, 4.836, 0.5243, 25),
    @(0.006925, 5.666, 0.1013, 10),
    @(0.02144, 6.043, 0.0873, 15),
    @(0.006183, 5.847, 0.5238, 15),
    @(0.02639, 6.084, 0.3582, 20)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
}